# Arquivo atualizado em 07/12/2023, 14:20:05.
#
# Column B held date values (serial numbers, formatted as
# "YYYY-MM-DD HH:MM:SS" via a custom number format / style) for the
# "Ano" (year) column. This converts those 30 cells (rows 2-31) into
# plain text cells holding the literal string "dd/mm/yyyy" (e.g.
# "01/01/2013"), and drops the now-unused date style from the cells so
# they fall back to the sheet's default (unstyled) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of the date serial that used to live in each row to the literal
# text that replaces it (they are all January 1st of a given year).
$yearText = @{
    2  = "01/01/2013"
    3  = "01/01/2014"
    4  = "01/01/2015"
    5  = "01/01/2016"
    6  = "01/01/2017"
    7  = "01/01/2018"
    8  = "01/01/2019"
    9  = "01/01/2020"
    10 = "01/01/2021"
    11 = "01/01/2022"
    12 = "01/01/2013"
    13 = "01/01/2014"
    14 = "01/01/2015"
    15 = "01/01/2016"
    16 = "01/01/2017"
    17 = "01/01/2018"
    18 = "01/01/2019"
    19 = "01/01/2020"
    20 = "01/01/2021"
    21 = "01/01/2022"
    22 = "01/01/2013"
    23 = "01/01/2014"
    24 = "01/01/2015"
    25 = "01/01/2016"
    26 = "01/01/2017"
    27 = "01/01/2018"
    28 = "01/01/2019"
    29 = "01/01/2020"
    30 = "01/01/2021"
    31 = "01/01/2022"
}

$dataRange = $ws.Range("B2:B31")

# Force the range to a text format first so the "dd/mm/yyyy"-looking
# strings we are about to write are stored verbatim instead of being
# re-interpreted as dates by Excel's usual autoconvert-on-entry logic.
$dataRange.NumberFormat = "@"

foreach ($row in 2..31) {
    $ws.Range("B$row").Value = $yearText[$row]
}

# The cells no longer need the custom date style that used to live
# here - reset them back to the workbook's default "Normal" style.
$dataRange.Style = "Normal"
